$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row before the current row 505. This shifts the
# existing rows 505-553 down to 506-554 (all of their contents/styles move
# with them), and leaves row 505 ready to receive the new record's values.
$ws.Rows(505).Insert()

$ws.Cells.Item(505, 1).Value2 = 4
$ws.Cells.Item(505, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(505, 3).Value2 = "Los Lagos"
$ws.Cells.Item(505, 4).Value2 = 44769
$ws.Cells.Item(505, 5).Value2 = 10
$ws.Cells.Item(505, 6).Value2 = "Fruta"
$ws.Cells.Item(505, 7).Value2 = 100108
$ws.Cells.Item(505, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(505, 9).Value2 = 100108006
$ws.Cells.Item(505, 10).Value2 = "Plátano"
$ws.Cells.Item(505, 11).Value2 = "Sin especificar"
$ws.Cells.Item(505, 12).Value2 = "Pintón"
$ws.Cells.Item(505, 13).Value2 = 300
$ws.Cells.Item(505, 14).Value2 = 32000
$ws.Cells.Item(505, 15).Value2 = 32000
$ws.Cells.Item(505, 16).Value2 = 32000
$ws.Cells.Item(505, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(505, 18).Value2 = "Ecuador"
$ws.Cells.Item(505, 19).Value2 = 1600
$ws.Cells.Item(505, 20).Value2 = 20
